# Auto-update draw results: append the 2025-11-13 Pick 4 draw as a new
# last row in the "Results" sheet (mirrors the daily scraper's commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the first empty row right after the current data (row 57 -> 58).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row   # xlUp = -4162
$newRow = $lastRow + 1

$date        = "2025-11-13"
$game        = "Pick 4"
$phase       = "251113"
$result      = "8-8-5-8"
$insertedAt  = "2025-11-13T21:41:09.801+04:00"

# Columns A, C and E hold values that *look* like dates/numbers
# ("2025-11-13", "251113", an ISO timestamp) but must be stored as plain
# text, exactly like every other row in the sheet. Pre-format those three
# cells as Text ("@") so Excel doesn't silently coerce them into a date
# serial / number on assignment, then drop the formatting again so the
# new row doesn't end up with a different cell style than the rest of
# the table (only the stored value type - text - needs to stick).
$dateCell       = $ws.Cells.Item($newRow, 1)
$phaseCell      = $ws.Cells.Item($newRow, 3)
$insertedCell   = $ws.Cells.Item($newRow, 5)

$dateCell.NumberFormat = "@"
$phaseCell.NumberFormat = "@"
$insertedCell.NumberFormat = "@"

$dateCell.Value = $date
$ws.Cells.Item($newRow, 2).Value = $game
$phaseCell.Value = $phase
$ws.Cells.Item($newRow, 4).Value = $result
$insertedCell.Value = $insertedAt

$dateCell.ClearFormats()
$phaseCell.ClearFormats()
$insertedCell.ClearFormats()
